# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled run).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells on this
# sheet (t="inlineStr"), so every new value below must land back in the cell
# as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values, keyed by cell address,
# exactly as they appear in the refreshed coinranking.com snapshot.
$updates = [ordered]@{
    'D2'  = '28.406.07';   'E2'  = '  +4.36%  '
    'D3'  = '1.724.92';    'E3'  = '  +2.19%  '
                           'E4'  = '  -0.05%  '
    'D5'  = '218.53';      'E5'  = '  +1.25%  '
                           'E6'  = '  +0.32%  '
    'D7'  = '1.00';        'E7'  = '  -0.15%  '
    'D8'  = '23.92';       'E8'  = '  +3.10%  '
                           'E9'  = '  +1.89%  '
    'D10' = '0.0634';      'E10' = '  +0.94%  '
    'D11' = '0.0892';      'E11' = '  +0.29%  '
    'D12' = '1.972.49';    'E12' = '  +2.44%  '
    'D13' = '1.729.39';    'E13' = '  +2.53%  '
    'D14' = '4.23';        'E14' = '  +0.54%  '
    'D15' = '0.563';       'E15' = '  +1.59%  '
    'D16' = '67.56';       'E16' = '  +0.55%  '
    'D17' = '28.363.90'
    'D18' = '246.79';      'E18' = '  +4.08%  '
    'D19' = '0.0₃0750';    'E19' = '  +0.51%  '
    'D20' = '7.87';        'E20' = '  -2.82%  '
                           'E21' = '  -0.24%  '
    'D22' = '4.60';        'E22' = '  +0.30%  '
    'D23' = '9.63';        'E23' = '  -0.39%  '
                           'E24' = '  -1.38%  '
    'D25' = '149.20';      'E25' = '  +1.30%  '
    'D26' = '7.42';        'E26' = '  +1.43%  '
    'D27' = '16.60';       'E27' = '  +0.88%  '
                           'E28' = '  +0.26%  '
                           'E29' = '  -0.10%  '
                           'E30' = '  +2.48%  '
                           'E31' = '  +3.02%  '
                           'E32' = '  +0.51%  '
    'D33' = '1.485.20';    'E33' = '  -4.47%  '
    'D34' = '3.24';        'E34' = '  -0.25%  '
                           'E35' = '  -2.12%  '
    'D36' = '0.976';       'E36' = '  +2.69%  '
                           'E37' = '  +0.43%  '
    'D38' = '0.599';       'E38' = '  -1.00%  '
                           'E39' = '  +1.56%  '
                           'E40' = '  -0.30%  '
    'D41' = '69.60';       'E41' = '  +0.33%  '
    'D42' = '1.00'
    'D43' = '5.65';        'E43' = '  -1.85%  '
    'D44' = '1.874.84';    'E44' = '  +1.91%  '
                           'E45' = '  +0.90%  '
    'D46' = '0.807';       'E46' = '  +1.96%  '
                           'E47' = '  +6.08%  '
    'D48' = '90.34';       'E48' = '  -0.75%  '
                           'E49' = '  +4.17%  '
    'D50' = '8.12';        'E50' = '  -2.58%  '
                           'E51' = '  -1.29%  '
}

# Several of the new Price figures (e.g. "1.00", "23.92") are digit strings
# that Excel's Range.Value setter would otherwise silently coerce into real
# numbers. Force Text format on just the column-D cells we are about to touch
# so the literal string round-trips unchanged, then drop the custom format
# again afterwards so the cell's style reverts to its original (unstyled)
# state instead of leaving a stray "@" number format behind.
$dAddresses = $updates.Keys | Where-Object { $_.StartsWith('D') }
foreach ($addr in $dAddresses) {
    $ws.Range($addr).NumberFormat = '@'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

foreach ($addr in $dAddresses) {
    $ws.Range($addr).Style = 'Normal'
}
